$d = $word.ActiveDocument

# 1. Remove "专家" from "IT技术专家：..." -> "IT技术：..."
$d.Content.Find.Execute(
    "IT技术专家：", $true, $false, $false, $false, $false,
    $true, 1, $false, "IT技术：", 2) | Out-Null

# 2. Re-home the "_GoBack" bookmark: remove it from its old spot (the lone
#    paragraph near the end of the document) and re-add it right before
#    "课程展示" inside the sentence we just edited. Word auto-splits the
#    run at the insertion point, which is what produces the extra <w:r>
#    boundaries seen in the diff.
$target = $d.Content
$target.Find.Execute("课程展示", $true) | Out-Null
$pos = $target.Start
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# 3. Split "IT" | "技术" | "：快速架构...变化的" into separate runs, matching
#    the diff, by briefly bookmarking the two internal boundaries (which
#    forces Word to break the run there) and then deleting those helper
#    bookmarks again - the run split persists.
$itRange = $d.Content
$itRange.Find.Execute("IT", $true) | Out-Null
$afterIT = $itRange.End

$jishuRange = $d.Range($afterIT, $d.Content.End)
$jishuRange.Find.Execute("技术", $true) | Out-Null
$afterJishu = $jishuRange.End

$d.Bookmarks.Add("zzTmpSplit1", $d.Range($afterIT, $afterIT)) | Out-Null
$d.Bookmarks.Add("zzTmpSplit2", $d.Range($afterJishu, $afterJishu)) | Out-Null

$d.Bookmarks.Item("zzTmpSplit1").Delete()
$d.Bookmarks.Item("zzTmpSplit2").Delete()
